$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column B (dates) from 2023-07-18 (serial 45125) to 2023-07-19 (serial 45126)
# for rows 2 through 61, preserving the existing date number format/style.
$ws.Range("B2:B61").Value = 45126
